# Arkteos reg3 decoder workbook update
# - add "Groupe Frigo" entries (fan speed / DC voltage / frigo + regulation errors)
# - re-order the R/S/T "export" helper columns (Python Dict Decoder now comes first)
# - re-categorise the outside-temperature row under "Groupe Frigo"
# (thanks @joel-bourquard)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Values")

# ---------------------------------------------------------------------------
# 1) Re-point the R/S/T header labels: R=Python Dict Decoder, S=OpenHAB Thing
#    Channels, T=OpenHAB items (was R=Thing Channels, S=items, T=Dict Decoder)
# ---------------------------------------------------------------------------
$ws.Range("R1").Value = "Python Dict Decoder"
$ws.Range("S1").Value = "OpenHAB Thing Channels"
$ws.Range("T1").Value = "OpenHAB items"

# ---------------------------------------------------------------------------
# 2) Rotate the R/S/T formulas for existing data rows 2-20:
#       new R (dict decoder)      <- old T formula
#       new S (thing channel)     <- old R formula
#       new T (openhab item)      <- old S formula
#    Row 21-23 extend the dict-decoder column further down.
# ---------------------------------------------------------------------------
$dictTmpl  = '=IF(H{0}="Oui",CONCAT("{{ ''stream'' : ",I{0},", ''name'' : ''",J{0},"'' ,''descr'' : ''",K{0},"'', ''byte1'': ",L{0},", ''weight1'': ",M{0},", ''byte2'': ",N{0},", ''weight2'': ",O{0},", ''divider'': ",P{0}," }},"),"")'
$thingTmpl = '=IF(H{0}="Oui",CONCAT("        Type number : ",J{0}," """,K{0},"""  [stateTopic=""arkteos/reg3/",J{0},"""]"),"")'
$itemTmpl  = '=IF(H{0}="Oui",CONCAT("        Number Arkteosreg3_",J{0}," """,K{0},""" {{channel=""mqtt:topic:arkteos-reg3:",J{0},"""}}"),"")'

for ($r = 2; $r -le 20; $r++) {
    $ws.Range("S$r").Formula = $thingTmpl -f $r
    $ws.Range("T$r").Formula = $itemTmpl  -f $r
}

# ---------------------------------------------------------------------------
# 3) "Température extérieure" now belongs to the "Groupe Frigo" ensemble
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = "Groupe Frigo"

# ---------------------------------------------------------------------------
# 4) New rows: fan speed + DC voltage (Groupe Frigo), regulation / frigo
#    active-error bytes.
# ---------------------------------------------------------------------------
$ws.Range("A21").Value = "Groupe Frigo"
$ws.Range("B21").Value = "Vitesse ventilateur (450 tr/min)"
$ws.Range("H21").Value = "Oui"
$ws.Range("I21").Value = 163
$ws.Range("J21").Value = "fan_speed_evaporator_1"
$ws.Range("K21").Value = "Vitesse ventalisateur groupe frigo 1"
$ws.Range("L21").Value = 56
$ws.Range("M21").Value = 1
$ws.Range("N21").Value = 57
$ws.Range("O21").Value = 256
$ws.Range("P21").Value = 1

$ws.Range("A22").Value = "Groupe Frigo"
$ws.Range("B22").Value = "Voltage DC"
$ws.Range("C22").Value = 163
$ws.Range("H22").Value = "Oui"
$ws.Range("I22").Value = 163
$ws.Range("J22").Value = "dc_voltage"
$ws.Range("K22").Value = "Voltage DC groupe frigo 1"
$ws.Range("L22").Value = 62
$ws.Range("M22").Value = 1
$ws.Range("N22").Value = 63
$ws.Range("O22").Value = 256
$ws.Range("P22").Value = 1

$ws.Range("A23").Value = "Général PAC"
$ws.Range("B23").Value = "Erreur active n°1  (régulation)"
$ws.Range("C23").Value = 227
$ws.Range("H23").Value = "Différemment"
$ws.Range("I23").Value = 227
$ws.Range("J23").Value = "active_error_reg"
$ws.Range("K23").Value = "Erreur régulation"
$ws.Range("L23").Value = 30
$ws.Range("M23").Value = 1
$ws.Range("N23").Value = 31
$ws.Range("O23").Value = 256
$ws.Range("P23").Value = 1
$ws.Range("Q23").Value = "& 0x0F"

$ws.Range("A24").Value = "Groupe Frigo"
$ws.Range("B24").Value = "Erreur active n°1 (frigo)"
$ws.Range("C24").Value = 163
$ws.Range("H24").Value = "Différemment"
$ws.Range("I24").Value = 163
$ws.Range("K24").Value = "Erreur groupe frigo"
$ws.Range("L24").Value = 12
$ws.Range("M24").Value = 1
$ws.Range("N24").Value = 13
$ws.Range("O24").Value = 256
$ws.Range("P24").Value = 1
$ws.Range("Q24").Value = "& 0x0F"

for ($r = 21; $r -le 23; $r++) {
    $ws.Range("R$r").Formula = $dictTmpl -f $r
}

# ---------------------------------------------------------------------------
# 5) The remaining "candidate" rows shift down two slots (21/22 used to be
#    blank, the fridge error rows now take 23/24) - rewrite them verbatim.
# ---------------------------------------------------------------------------
$ws.Range("A25").Value = "Chauffage"
$ws.Range("B25").Value = "Débit eau ?"

$ws.Range("A26").Value = "Chauffage"
$ws.Range("B26").Value = "% circulateur ?"

$ws.Range("A27").Value = "Chauffage"
$ws.Range("B27").Value = "Etat résistance  ?"

$ws.Range("A28").ClearContents()
$ws.Range("B28").Value = "Statut PAC (chaud / froid)"

$ws.Range("A29").Value = "Général PAC"
$ws.Range("B29").Value = "PW produit instant"

$ws.Range("A30").Value = "Général PAC"
$ws.Range("B30").Value = "PW conso instant"
$ws.Range("E30").Value = "227!156"

# ---------------------------------------------------------------------------
# 6) Swap the R/S/T column widths to match their new contents
# ---------------------------------------------------------------------------
$ws.Columns.Item(18).ColumnWidth = 70.18
$ws.Columns.Item(19).ColumnWidth = 32.73
$ws.Columns.Item(20).ColumnWidth = 127.27

# ---------------------------------------------------------------------------
# 7) Selection / scroll position as left by the editor
# ---------------------------------------------------------------------------
$ws.Range("O24").Select()
